$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text in the source data (e.g. "606.05",
# "69.363.68"), not real numbers. Force the whole column to Text format
# before writing so numeric-looking strings ("603.19", "1.00", ...) are not
# silently reinterpreted as numbers; restore the style afterwards so no
# visible per-cell formatting change is left behind.
$dcol = $ws.Range("D2:D51")
$dcol.NumberFormat = "@"

# --- Rows 38/39: EthereumClassic and Monero swap rank positions ---
# Row 38 becomes Monero (was EthereumClassic), Row 39 becomes EthereumClassic (was Monero)
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "162.69"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "20.06"
$ws.Range("E39").Value = "  +1.59%  "

# --- Rows 45/46: USDe and BabyDogeCoin swap rank positions ---
# Row 45 becomes USDe (was BabyDogeCoin), Row 46 becomes BabyDogeCoin (was USDe)
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0317"
$ws.Range("E46").Value = "  -2.76%  "

# --- Remaining price/volume refresh updates ---
$ws.Range("D2").Value = "69.080.26"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "2.742.12"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "603.19"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "169.21"
$ws.Range("E6").Value = "  +6.45%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "2.740.27"
$ws.Range("E9").Value = "  +3.73%  "
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "0.368"
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "28.95"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").Value = "3.242.61"
$ws.Range("E15").Value = "  +3.86%  "
$ws.Range("D16").Value = "0.0000191"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "68.988.68"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "2.744.96"
$ws.Range("E18").Value = "  +4.08%  "
$ws.Range("D19").Value = "11.89"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").Value = "373.65"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "7.75"
$ws.Range("E21").Value = "  +5.63%  "
$ws.Range("D22").Value = "4.56"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").Value = "5.00"
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("D25").Value = "74.06"
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "2.881.85"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").Value = "0.0000107"
$ws.Range("E29").Value = "  +3.41%  "
$ws.Range("D30").Value = "599.06"
$ws.Range("E30").Value = "  +6.31%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "8.37"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").Value = "1.46"
$ws.Range("E33").Value = "  +5.23%  "
$ws.Range("D34").Value = "1.98"
$ws.Range("E34").Value = "  +6.93%  "
$ws.Range("D35").Value = "0.134"
$ws.Range("E35").Value = "  +4.74%  "
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D40").Value = "0.384"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("D41").Value = "1.93"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("D42").Value = "5.54"
$ws.Range("E42").Value = "  +4.54%  "
$ws.Range("D43").Value = "2.71"
$ws.Range("E43").Value = "  +4.07%  "
$ws.Range("D44").Value = "18.00"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D47").Value = "40.93"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "157.96"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "3.98"
$ws.Range("E49").Value = "  +6.55%  "
$ws.Range("E50").Value = "  +8.05%  "
$ws.Range("D51").Value = "0.609"
$ws.Range("E51").Value = "  +7.91%  "

# Restore the original (default/"Normal") style on column D now that all the
# text values are in place, so the cells carry no leftover number-format
# override.
$dcol.Style = "Normal"
